$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Tue Oct 08 21:00:34 BRT 2024"
$ws.Range("C3").Value = "entrega"
